# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new value for column F
$wsExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    4  = 13267
    7  = 217
    13 = 13233
    16 = 8853
    17 = 7931
    24 = 1007
    27 = 390
    29 = 113
    30 = 358
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    5  = 13267
    8  = 217
    14 = 13233
    17 = 8853
    18 = 7931
    25 = 1007
    30 = 390
    32 = 113
    33 = 358
}
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
